$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-16 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-17 Monday", 2) | Out-Null
$d.Content.Find.Execute("98÷8=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "25÷4=6, 1", 2) | Out-Null
$d.Content.Find.Execute("64÷3=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "21÷7=3, 0", 2) | Out-Null
$d.Content.Find.Execute("56÷5=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "98÷3=32, 2", 2) | Out-Null
$d.Content.Find.Execute("86÷2=43, 0", $true, $false, $false, $false, $false, $true, 1, $false, "29÷5=5, 4", 2) | Out-Null
$d.Content.Find.Execute("77÷6=12, 5", $true, $false, $false, $false, $false, $true, 1, $false, "66÷7=9, 3", 2) | Out-Null
$d.Content.Find.Execute("87÷9=9, 6", $true, $false, $false, $false, $false, $true, 1, $false, "97÷8=12, 1", 2) | Out-Null
$d.Content.Find.Execute("55÷6=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "75÷4=18, 3", 2) | Out-Null
$d.Content.Find.Execute("53÷4=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "40÷7=5, 5", 2) | Out-Null
$d.Content.Find.Execute("29÷4=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "84÷4=21, 0", 2) | Out-Null
$d.Content.Find.Execute("74÷9=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "42÷4=10, 2", 2) | Out-Null
$d.Content.Find.Execute("46÷4=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "77÷4=19, 1", 2) | Out-Null
$d.Content.Find.Execute("46÷6=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "97÷2=48, 1", 2) | Out-Null
$d.Content.Find.Execute("24÷8=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "68÷4=17, 0", 2) | Out-Null
$d.Content.Find.Execute("14÷7=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "71÷2=35, 1", 2) | Out-Null
$d.Content.Find.Execute("23÷8=2, 7", $true, $false, $false, $false, $false, $true, 1, $false, "28÷6=4, 4", 2) | Out-Null
$d.Content.Find.Execute("81÷6=13, 3", $true, $false, $false, $false, $false, $true, 1, $false, "76÷7=10, 6", 2) | Out-Null
$d.Content.Find.Execute("13÷4=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "65÷9=7, 2", 2) | Out-Null
$d.Content.Find.Execute("42÷9=4, 6", $true, $false, $false, $false, $false, $true, 1, $false, "22÷4=5, 2", 2) | Out-Null
$d.Content.Find.Execute("90÷6=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "70÷4=17, 2", 2) | Out-Null
$d.Content.Find.Execute("45÷8=5, 5", $true, $false, $false, $false, $false, $true, 1, $false, "14÷7=2, 0", 2) | Out-Null
$d.Content.Find.Execute("49÷9=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "29÷2=14, 1", 2) | Out-Null
$d.Content.Find.Execute("50÷9=5, 5", $true, $false, $false, $false, $false, $true, 1, $false, "48÷5=9, 3", 2) | Out-Null
$d.Content.Find.Execute("75÷3=25, 0", $true, $false, $false, $false, $false, $true, 1, $false, "62÷9=6, 8", 2) | Out-Null
$d.Content.Find.Execute("86÷5=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "79÷2=39, 1", 2) | Out-Null
$d.Content.Find.Execute("42÷5=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "26÷4=6, 2", 2) | Out-Null
